# Auto-generated edit script: applies the numeric cell updates described by the
# commit diff across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 200
$ws.Range("K29").Value = 600
$ws.Range("M29").Value = -319

$ws.Range("H58").Value = 715243.2
$ws.Range("J58").Value = 2857.1428
$ws.Range("L58").Value = 8571.428400000001
$ws.Range("N58").Value = -8871.428400000001

$ws.Range("H106").Value = 1910.2142
$ws.Range("I106").Value = 1631.1818
$ws.Range("K106").Value = 1631.1818
$ws.Range("M106").Value = -1000.1818

$ws.Range("H129").Value = 2154.8684
$ws.Range("I129").Value = 6750.875
$ws.Range("J129").Value = 929.26666
$ws.Range("K129").Value = 20252.625
$ws.Range("L129").Value = 2787.79998
$ws.Range("M129").Value = -15252.625
$ws.Range("N129").Value = -12787.79998

$ws.Range("H132").Value = 3475709.8
$ws.Range("I132").Value = 3574854.2
$ws.Range("J132").Value = 5658.5
$ws.Range("K132").Value = 10724562.6
$ws.Range("L132").Value = 16975.5
$ws.Range("M132").Value = -10722032.6
$ws.Range("N132").Value = -22035.5

$ws.Range("H137").Value = 1163.1708
$ws.Range("I137").Value = 1079.2307
$ws.Range("J137").Value = 2800
$ws.Range("K137").Value = 3237.6921
$ws.Range("L137").Value = 8400
$ws.Range("M137").Value = -687.6921000000002
$ws.Range("N137").Value = -13500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H61").Value = 1788.5193
$ws.Range("I61").Value = 929.96155
$ws.Range("J61").Value = 2647.077
$ws.Range("K61").Value = 929.96155
$ws.Range("L61").Value = 2647.077
$ws.Range("M61").Value = -717.96155
$ws.Range("N61").Value = -3071.077

$ws.Range("H74").Value = 763.55554
$ws.Range("I74").Value = 774.37036
$ws.Range("J74").Value = 731.1111
$ws.Range("K74").Value = 774.37036
$ws.Range("L74").Value = 731.1111
$ws.Range("M74").Value = 99.62963999999999
$ws.Range("N74").Value = -2479.1111

$ws.Range("H77").Value = 763.55554
$ws.Range("I77").Value = 774.37036
$ws.Range("J77").Value = 731.1111
$ws.Range("K77").Value = 3871.8518
$ws.Range("L77").Value = 3655.5555
$ws.Range("M77").Value = 496.1482000000001
$ws.Range("N77").Value = -12391.5555

$ws.Range("H110").Value = 55672824
$ws.Range("I110").Value = 83508760
$ws.Range("K110").Value = 83508760
$ws.Range("M110").Value = -83506715

$ws.Range("H131").Value = 34572
$ws.Range("J131").Value = 34572
$ws.Range("L131").Value = 34572
$ws.Range("N131").Value = -44652

$ws.Range("H132").Value = 2105.111
$ws.Range("I132").Value = 2029.1296
$ws.Range("J132").Value = 2561
$ws.Range("K132").Value = 6087.3888
$ws.Range("L132").Value = 7683
$ws.Range("M132").Value = -3557.3888
$ws.Range("N132").Value = -12743

$ws.Range("H136").Value = 1788.5193
$ws.Range("I136").Value = 929.96155
$ws.Range("J136").Value = 2647.077
$ws.Range("K136").Value = 2789.88465
$ws.Range("L136").Value = 7941.231000000001
$ws.Range("M136").Value = -239.88465
$ws.Range("N136").Value = -13041.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43984.082
$ws.Range("I20").Value = 61062.883
$ws.Range("J20").Value = 2507
$ws.Range("K20").Value = 61062.883
$ws.Range("L20").Value = 2507
$ws.Range("M20").Value = -60815.883
$ws.Range("N20").Value = -3001

$ws.Range("H94").Value = 676.1053000000001
$ws.Range("I94").Value = 502.36365
$ws.Range("J94").Value = 915
$ws.Range("K94").Value = 502.36365
$ws.Range("L94").Value = 915
$ws.Range("M94").Value = -51.36365000000001
$ws.Range("N94").Value = -1817

$ws.Range("H134").Value = 11002
$ws.Range("I134").Value = 11002.4
$ws.Range("J134").Value = 11000
$ws.Range("K134").Value = 33007.2
$ws.Range("L134").Value = 33000
$ws.Range("M134").Value = -30472.2
$ws.Range("N134").Value = -38070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 459.7143
$ws.Range("I22").Value = 478
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 478
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -128
$ws.Range("N22").Value = -1050

$ws.Range("H31").Value = 23730.625
$ws.Range("I31").Value = 1026.0588
$ws.Range("J31").Value = 49462.465
$ws.Range("K31").Value = 1026.0588
$ws.Range("L31").Value = 49462.465
$ws.Range("M31").Value = -731.0588
$ws.Range("N31").Value = -50052.465

$ws.Range("H34").Value = 23730.625
$ws.Range("I34").Value = 1026.0588
$ws.Range("J34").Value = 49462.465
$ws.Range("K34").Value = 1026.0588
$ws.Range("L34").Value = 49462.465
$ws.Range("M34").Value = -824.0588
$ws.Range("N34").Value = -49866.465

$ws.Range("H99").Value = 15997
$ws.Range("I99").Value = 5774.2
$ws.Range("J99").Value = 23299
$ws.Range("K99").Value = 5774.2
$ws.Range("L99").Value = 23299
$ws.Range("M99").Value = -4276.2
$ws.Range("N99").Value = -26295

$ws.Range("H107").Value = 846.5
$ws.Range("I107").Value = 1064.909
$ws.Range("J107").Value = 579.55554
$ws.Range("K107").Value = 1064.909
$ws.Range("L107").Value = 579.55554
$ws.Range("M107").Value = 855.0909999999999
$ws.Range("N107").Value = -4419.55554

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 15997
$ws.Range("I126").Value = 5774.2
$ws.Range("J126").Value = 23299
$ws.Range("K126").Value = 17322.6
$ws.Range("L126").Value = 69897
$ws.Range("M126").Value = -14852.6
$ws.Range("N126").Value = -74837

$ws.Range("H132").Value = 37503356
$ws.Range("I132").Value = 31253304
$ws.Range("J132").Value = 62503564
$ws.Range("K132").Value = 93759912
$ws.Range("L132").Value = 187510692
$ws.Range("M132").Value = -93757382
$ws.Range("N132").Value = -187515752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1838

$ws.Range("H34").Value = 687.5714
$ws.Range("I34").Value = 362.16666
$ws.Range("J34").Value = 931.625
$ws.Range("K34").Value = 1086.49998
$ws.Range("L34").Value = 2794.875
$ws.Range("M34").Value = -1002.49998
$ws.Range("N34").Value = -2962.875

$ws.Range("H39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H44").Value = 778.8889
$ws.Range("I44").Value = 251
$ws.Range("K44").Value = 753
$ws.Range("M44").Value = -355

$ws.Range("H111").Value = 1113.5714
$ws.Range("I111").Value = 1113.5714
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3340.7142
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -273.7142000000003
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 85476.414
$ws.Range("I112").Value = 334409
$ws.Range("J112").Value = 2498.889
$ws.Range("K112").Value = 1003227
$ws.Range("L112").Value = 7496.667
$ws.Range("M112").Value = -1002119
$ws.Range("N112").Value = -9712.667000000001

$ws.Range("H131").Value = 11076.085
$ws.Range("J131").Value = 11562.143
$ws.Range("L131").Value = 34686.429
$ws.Range("N131").Value = -44766.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2648816.2
$ws.Range("I11").Value = 2793562
$ws.Range("J11").Value = 550000
$ws.Range("K11").Value = 2793562
$ws.Range("L11").Value = 550000
$ws.Range("M11").Value = -2793423
$ws.Range("N11").Value = -550278

$ws.Range("H70").Value = 51745.906
$ws.Range("I70").Value = 90562.56
$ws.Range("J70").Value = 4757.316
$ws.Range("K70").Value = 90562.56
$ws.Range("L70").Value = 4757.316
$ws.Range("M70").Value = -90292.56
$ws.Range("N70").Value = -5297.316

$ws.Range("H73").Value = 51745.906
$ws.Range("I73").Value = 90562.56
$ws.Range("J73").Value = 4757.316
$ws.Range("K73").Value = 90562.56
$ws.Range("L73").Value = 4757.316
$ws.Range("M73").Value = -89626.56
$ws.Range("N73").Value = -6629.316

$ws.Range("H107").Value = 561519.1
$ws.Range("I107").Value = 390.15384
$ws.Range("J107").Value = 2020454.4
$ws.Range("K107").Value = 390.15384
$ws.Range("L107").Value = 2020454.4
$ws.Range("M107").Value = 1529.84616
$ws.Range("N107").Value = -2024294.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 818.1515000000001
$ws.Range("I22").Value = 772.63635
$ws.Range("J22").Value = 840.9091
$ws.Range("K22").Value = 772.63635
$ws.Range("L22").Value = 840.9091
$ws.Range("M22").Value = -477.63635
$ws.Range("N22").Value = -1430.9091

$ws.Range("H27").Value = 818.1515000000001
$ws.Range("I27").Value = 772.63635
$ws.Range("J27").Value = 840.9091
$ws.Range("K27").Value = 772.63635
$ws.Range("L27").Value = 840.9091
$ws.Range("M27").Value = -665.63635
$ws.Range("N27").Value = -1054.9091

$ws.Range("H108").Value = 26875.334
$ws.Range("J108").Value = 26875.334
$ws.Range("L108").Value = 26875.334
$ws.Range("N108").Value = -34555.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2173.0227
$ws.Range("I132").Value = 2179.8206
$ws.Range("J132").Value = 2120
$ws.Range("K132").Value = 6539.4618
$ws.Range("L132").Value = 6360
$ws.Range("M132").Value = -4009.4618
$ws.Range("N132").Value = -11420

$ws.Range("H135").Value = 43287.445
$ws.Range("J135").Value = 43287.445
$ws.Range("L135").Value = 43287.445
$ws.Range("N135").Value = -53427.445
